# Updates the cryptos list price/volume columns (and swaps the Monero / LidoDAOToken
# row order) per the Tue Mar 21 17:45:49 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values that look numeric are written through a
# text-formatted cell (then restored to the default "Normal" style) so Excel
# stores the exact original string (trailing zeros, thousand-dot grouping, etc.)
# instead of silently coercing it to a Double and mangling the formatting.
$updates = @(
    @{Cell='D2'; Value='28.247.56'}
    @{Cell='E2'; Value='  +1.52%  '}
    @{Cell='D3'; Value='1.802.52'}
    @{Cell='E3'; Value='  +2.11%  '}
    @{Cell='D4'; Value='1.004'}
    @{Cell='E4'; Value='  -0.15%  '}
    @{Cell='D5'; Value='338.61'}
    @{Cell='E5'; Value='  -0.03%  '}
    @{Cell='D6'; Value='1.001'}
    @{Cell='E6'; Value='  +0.00%  '}
    @{Cell='D7'; Value='0.4498'}
    @{Cell='E7'; Value='  +19.32%  '}
    @{Cell='D8'; Value='0.3533'}
    @{Cell='E8'; Value='  +5.27%  '}
    @{Cell='D9'; Value='45.50'}
    @{Cell='E9'; Value='  -0.26%  '}
    @{Cell='D10'; Value='1.145'}
    @{Cell='E10'; Value='  +0.91%  '}
    @{Cell='D11'; Value='0.07560'}
    @{Cell='E11'; Value='  +4.15%  '}
    @{Cell='D12'; Value='22.67'}
    @{Cell='E12'; Value='  -1.18%  '}
    @{Cell='D13'; Value='1.004'}
    @{Cell='E13'; Value='  +0.06%  '}
    @{Cell='D14'; Value='6.233'}
    @{Cell='E14'; Value='  -0.22%  '}
    @{Cell='D15'; Value='7.238'}
    @{Cell='E15'; Value='  +0.22%  '}
    @{Cell='D16'; Value='1.800.00'}
    @{Cell='E16'; Value='  +1.85%  '}
    @{Cell='D17'; Value='0.00001089'}
    @{Cell='E17'; Value='  +3.21%  '}
    @{Cell='E18'; Value='  +1.12%  '}
    @{Cell='D19'; Value='81.41'}
    @{Cell='E19'; Value='  +0.33%  '}
    @{Cell='D20'; Value='0.9987'}
    @{Cell='E20'; Value='  -0.32%  '}
    @{Cell='D21'; Value='17.18'}
    @{Cell='E21'; Value='  +0.24%  '}
    @{Cell='D22'; Value='6.380'}
    @{Cell='E22'; Value='  +0.85%  '}
    @{Cell='D23'; Value='28.218.94'}
    @{Cell='E23'; Value='  +1.35%  '}
    @{Cell='D24'; Value='11.93'}
    @{Cell='E24'; Value='  +1.12%  '}
    @{Cell='D25'; Value='2.401'}
    @{Cell='E25'; Value='  +1.04%  '}
    @{Cell='D26'; Value='20.60'}
    @{Cell='E26'; Value='  +3.01%  '}
    @{Cell='B27'; Value='Monero'}
    @{Cell='C27'; Value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'}
    @{Cell='D27'; Value='155.10'}
    @{Cell='E27'; Value='  +1.12%  '}
    @{Cell='B28'; Value='LidoDAOToken'}
    @{Cell='C28'; Value='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'}
    @{Cell='D28'; Value='2.412'}
    @{Cell='E28'; Value='  +2.37%  '}
    @{Cell='D29'; Value='2.009.87'}
    @{Cell='E29'; Value='  +2.17%  '}
    @{Cell='D30'; Value='1.285'}
    @{Cell='E30'; Value='  -13.73%  '}
    @{Cell='D31'; Value='132.59'}
    @{Cell='E31'; Value='  +0.15%  '}
    @{Cell='D32'; Value='4.076'}
    @{Cell='D33'; Value='5.917'}
    @{Cell='E33'; Value='  +0.06%  '}
    @{Cell='D34'; Value='0.09442'}
    @{Cell='E34'; Value='  +8.02%  '}
    @{Cell='D35'; Value='0.02376'}
    @{Cell='E35'; Value='  +0.65%  '}
    @{Cell='D36'; Value='12.16'}
    @{Cell='E36'; Value='  -2.14%  '}
    @{Cell='D37'; Value='0.6696'}
    @{Cell='E37'; Value='  -0.15%  '}
    @{Cell='D38'; Value='0.06266'}
    @{Cell='E38'; Value='  -0.02%  '}
    @{Cell='D39'; Value='0.2158'}
    @{Cell='E39'; Value='  +1.91%  '}
    @{Cell='D40'; Value='5.180'}
    @{Cell='E40'; Value='  -0.48%  '}
    @{Cell='D41'; Value='1.485'}
    @{Cell='E41'; Value='  +0.88%  '}
    @{Cell='D42'; Value='1.213'}
    @{Cell='E42'; Value='  -0.91%  '}
    @{Cell='D43'; Value='8.157'}
    @{Cell='E43'; Value='  +0.92%  '}
    @{Cell='D44'; Value='0.9983'}
    @{Cell='E44'; Value='  -0.30%  '}
    @{Cell='D45'; Value='13.95'}
    @{Cell='E45'; Value='  +0.87%  '}
    @{Cell='D46'; Value='3.862'}
    @{Cell='E46'; Value='  +0.54%  '}
    @{Cell='D47'; Value='0.6096'}
    @{Cell='E47'; Value='  -0.46%  '}
    @{Cell='D48'; Value='129.13'}
    @{Cell='E48'; Value='  -1.77%  '}
    @{Cell='D49'; Value='2.028'}
    @{Cell='E49'; Value='  +0.07%  '}
    @{Cell='D50'; Value='0.07101'}
    @{Cell='E50'; Value='  -2.46%  '}
    @{Cell='D51'; Value='1.163'}
    @{Cell='E51'; Value='  -1.75%  '}
)

$applied = 0
foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $val = $u.Value
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        # numeric-looking text (price column) -> force text storage, then drop the
        # temporary number format so the cell's style matches the untouched cells.
        $c.NumberFormat = '@'
        $c.Value = $val
        $c.Style = 'Normal'
    } else {
        $c.Value = $val
    }
    $applied++
}

Write-Output "Applied $applied cell updates"
